# Major improvements to preventing back-to-back scheduling issues.
#
# - Abbreviations for "Major" and "Minor" divisions are set to match their
#   full names (instead of "Maj"/"Min").
# - maxLateGames (column E) reduced to 10 for Major/Minor/PeeWee.
# - maxGames (column F) reduced from 50 to 16 for every division.
# - daysBetween (column H) increased from 1 to 2 for every division, to
#   prevent teams from being scheduled for back-to-back games.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tablib Dataset")

# Column C: abbreviation -- Major/Minor rows now spell out the full name.
$ws.Range("C2").Value = "Major"
$ws.Range("C3").Value = "Minor"

# Column E: maxLateGames
$ws.Range("E2").Value = 10
$ws.Range("E3").Value = 10
$ws.Range("E4").Value = 10

# Column F: maxGames -- every division drops from 50 to 16.
$ws.Range("F2:F8").Value = 16

# Column H: daysBetween -- every division bumps from 1 to 2.
$ws.Range("H2:H8").Value = 2

# Selection, as left by the editor.
$ws.Range("M13").Select() | Out-Null
